# Refactor api response formatting
# Appends a new row (51) to each of the 4 worksheets, mirroring the
# structure of the existing row 50 (same layout: date, two text/byte
# strings, a checksum string, a byte-count string, and three numeric
# counters).

$wb = $excel.ActiveWorkbook

$rowsData = @{
    "MID_LFT_#1" = @{
        A = 45837.46297453704
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
        D = "0x01,0x6C"
        E = "0x07"
        F = 400
        G = [double]"5.68631262647113e+23"
        H = 364
        I = 7
    }
    "MID_LFT_#2" = @{
        A = 45837.46297453704
        B = "0x01,0x7c"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
        D = "0x01,0x60"
        E = "0x19"
        F = 380
        G = [double]"5.68432987514711e+23"
        H = 352
        I = 25
    }
    "MID_PLT_#1" = @{
        A = 45837.46297453704
        B = "0x00,0x6e"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
        D = "0x00,0x69"
        E = "0x15"
        F = 110
        G = [double]"5.68631262647113e+23"
        H = 105
        I = 15
    }
    "MID_PLT_#2" = @{
        A = 45837.46297453704
        B = "0x00,0x82"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
        D = "0x00,0x7E"
        E = "0x9"
        F = 130
        G = [double]"5.68631262647113e+23"
        H = 126
        I = 9
    }
}

foreach ($ws in $wb.Worksheets) {
    $name = $ws.Name
    if (-not $rowsData.ContainsKey($name)) { continue }
    $data = $rowsData[$name]

    $newRow = 51

    # Column A: date/time value, same number format as the row above it.
    $cellA = $ws.Cells.Item($newRow, 1)
    $cellA.NumberFormat = $ws.Cells.Item($newRow - 1, 1).NumberFormat
    $cellA.Value = $data.A

    # Columns B-E: text-ish byte strings (stored as text).
    $ws.Cells.Item($newRow, 2).Value = $data.B
    $ws.Cells.Item($newRow, 3).Value = $data.C
    $ws.Cells.Item($newRow, 4).Value = $data.D
    $ws.Cells.Item($newRow, 5).Value = $data.E

    # Columns F-I: plain numbers.
    $ws.Cells.Item($newRow, 6).Value = $data.F
    $ws.Cells.Item($newRow, 7).Value = $data.G
    $ws.Cells.Item($newRow, 8).Value = $data.H
    $ws.Cells.Item($newRow, 9).Value = $data.I
}
